# Weekly update: a new price-report row for the current week is inserted
# at row 167 (just before the existing 2023-09-10 / 45033 entry), pushing
# every subsequent row down by one. The former last row (299) becomes the
# new last row (300).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 167:299 down to 168:300, inheriting row 167's former formatting.
$ws.Rows.Item(167).Insert()

# Populate the newly-opened row 167 with this week's data.
$ws.Cells.Item(167, 1).Value = 5
$ws.Cells.Item(167, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(167, 3).Value = "Maule"
$ws.Cells.Item(167, 4).Value = 45072
$ws.Cells.Item(167, 5).Value = 7
$ws.Cells.Item(167, 6).Value = 100112017
$ws.Cells.Item(167, 7).Value = "Apio"
$ws.Cells.Item(167, 8).Value = "Americana (o)"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 700
$ws.Cells.Item(167, 11).Value = 6000
$ws.Cells.Item(167, 12).Value = 6000
$ws.Cells.Item(167, 13).Value = 6000
$ws.Cells.Item(167, 14).Value = "`$/docena de matas"
$ws.Cells.Item(167, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(167, 16).Value = 1000
$ws.Cells.Item(167, 17).Value = 6
$ws.Cells.Item(167, 18).Value = "Hortaliza"
